# Updated to 29 February
# - refresh the archived-source timestamp for the 27.02.2024 row (M143)
# - correct the running totals for 27.02.2024 (C143, F143)
# - append two new daily rows for 28.02.2024 and 29.02.2024

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing row 143 (27.02.2024) -------------------------------
$ws.Range("C143").Value = 29878
$ws.Range("F143").Value = 70215
$ws.Range("M143").Value = "https://web.archive.org/web/20240227162215/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- New row 144 (28.02.2024) -------------------------------------------
$ws.Range("A144").Value = "28.02.2024"
$ws.Range("B144").Value = "28.02.2024"
$ws.Range("C144").Value = 29954
$ws.Range("D144").Value = 12300
$ws.Range("E144").Value = 8400
$ws.Range("F144").Value = 70325
$ws.Range("G144").Value = 8663
$ws.Range("H144").Value = 6327
$ws.Range("I144").Value = 7000
$ws.Range("J144").Value = 411
$ws.Range("K144").Value = 108
$ws.Range("L144").Value = 4600
$ws.Range("M144").Value = "https://web.archive.org/web/20240228183950/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- New row 145 (29.02.2024) -------------------------------------------
$ws.Range("A145").Value = "29.02.2024"
$ws.Range("B145").Value = "28.02.2024"
$ws.Range("C145").Value = 30035
$ws.Range("D145").Value = 12300
$ws.Range("E145").Value = 8400
$ws.Range("F145").Value = 70457
$ws.Range("G145").Value = 8663
$ws.Range("H145").Value = 6327
$ws.Range("I145").Value = 8000
$ws.Range("J145").Value = 414
$ws.Range("K145").Value = 108
$ws.Range("L145").Value = 4600
$ws.Range("M145").Value = "https://web.archive.org/web/20240229162901/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Best-effort view-state refresh (selection / split) -----------------
$ws.Range("I148").Select()
$win = $excel.ActiveWindow
$win.SplitRow = 121
